# actualizacion Vo.Bo. UPP 4T 2020 SIPOT
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the stray formatted cell left over in row 11 (G11) ---
$ws.Rows("11:11").Delete()

# --- Remove the now-unused column I (merges collapse automatically) ---
$ws.Columns("I:I").Delete()

# --- Update the reporting-period dates in row 8 for the new (4T 2020) period ---
$ws.Range("B8").Value = 44105   # 10/01/2020
$ws.Range("C8").Value = 44196   # 12/31/2020
$ws.Range("F8").Value = 44206   # 01/10/2021
$ws.Range("G8").Value = 44206   # 01/10/2021

# --- Re-apply a full border box to the cells that lost their neighbour when
#     column I was deleted (they used to share a border with the merged I cell) ---
$ws.Range("H7").Borders.LineStyle = 1
$ws.Range("H8").Borders.LineStyle = 1
$ws.Range("H8").HorizontalAlignment = -4131   # xlLeft, matches the rest of row 8

# --- G3's description header switches from left to justified wrap alignment ---
$ws.Range("G3").HorizontalAlignment = -4130   # xlJustify

# --- Row heights grew slightly to fit the rewrapped text ---
$ws.Rows(3).RowHeight = 65.25
$ws.Rows(8).RowHeight = 54.75

# --- Column widths were tweaked (D narrower, G/H wider to fit the new layout) ---
$ws.Columns(4).ColumnWidth = 71.8
$ws.Columns(7).ColumnWidth = 21.65
$ws.Columns(8).ColumnWidth = 56.3

# --- Selection moved to the title cell block ---
$ws.Range("A3:C3").Select()

# --- Page was set up for printing on letter paper, portrait orientation ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
